$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.680.89"
$ws.Range("E2").Value = "  -0.27%  "

$ws.Range("D3").Value = "2.289.88"
$ws.Range("E3").Value = "  -1.25%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "103.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "270.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.624"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.41%  "

$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.607"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.67%  "

$ws.Range("E11").Value = "  -1.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.35%  "

$ws.Range("E13").Value = "  +1.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.88%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.857"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.83%  "

$ws.Range("D16").Value = "2.293.15"
$ws.Range("E16").Value = "  -1.28%  "

$ws.Range("D17").Value = "43.685.18"
$ws.Range("E17").Value = "  -0.20%  "

$ws.Range("E18").Value = "  +0.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "233.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.71%  "

$ws.Range("E23").Value = "  +14.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.39%  "

$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.51%  "

$ws.Range("E27").Value = "  -0.99%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "40.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.78%  "

$ws.Range("E29").Value = "  -2.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "177.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0899"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.77%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.47%  "

$ws.Range("E35").Value = "  -0.27%  "

$ws.Range("E36").Value = "  +0.05%  "

$ws.Range("E37").Value = "  -2.79%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.54"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.238"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.12%  "

$ws.Range("E40").Value = "  -0.87%  "

$ws.Range("E41").Value = "  +0.48%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.97%  "

$ws.Range("E44").Value = "  -1.92%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.51%  "

$ws.Range("E46").Value = "  -1.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.95%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.66%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.448"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.57%  "

$ws.Range("D51").Value = "2.513.90"
$ws.Range("E51").Value = "  -1.22%  "
